$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed "K" column (column G) values: strikeouts replaced with K-count based values
$kValues = @(0, 1, 0, 2, 0, 1, 0, 1, 1, 0, 2, 1, 2, 2, 2, 1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
